# Applies the LOT2055.xlsx edit:
#  - remove the last row (25) so the sheet shrinks to A1:C24
#  - rewrite every remaining cell to its final text
#  - fix up custom row heights / clear cells that no longer apply
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(25).Delete()

# Row 1
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Rows(1).EntireRow.AutoFit()

# Row 2
$ws.Range("B2").Value = 'LOT2055'
$ws.Range("C2").Value = 'LOT2055'
$ws.Rows(2).EntireRow.AutoFit()

# Row 3
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Estatística em Bioprocessos'
$ws.Range("C3").Value = ' Estatística em Bioprocessos'
$ws.Rows(3).EntireRow.AutoFit()

# Row 4
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Statistic for Bioprocesses'
$ws.Range("C4").Value = 'Statistic for Bioprocesses'
$ws.Rows(4).EntireRow.AutoFit()

# Row 5
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '2'
$ws.Range("C5").Value = '2'
$ws.Rows(5).EntireRow.AutoFit()

# Row 6
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '1'
$ws.Range("C6").Value = '1'
$ws.Rows(6).EntireRow.AutoFit()

# Row 7
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'
$ws.Rows(7).EntireRow.AutoFit()

# Row 8
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2019'
$ws.Range("C8").Value = '01/01/2019'
$ws.Rows(8).EntireRow.AutoFit()

# Row 9
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EB-9'
$ws.Range("C9").Value = 'EB-9'
$ws.Rows(9).EntireRow.AutoFit()

# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '5817181 - Valdeir Arantes'
$ws.Range("C10").Value = '5817181 - Valdeir Arantes'
$ws.Rows(10).RowHeight = 60

# Row 11
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'To familiarize the student with the basic concepts of applied statistics to study the influence of independent variables on dependent variables (responses) in bioprocesses. To Introduce to the student experimental design tools used to plan, perform full and fractional factorial designs, analyze the results, model the process based on the empirical data and suggest optimization conditions. To familiarize the student with a commercial software on the subject.'
$ws.Range("C11").Value = 'To familiarize the student with the basic concepts of applied statistics to study the influence of independent variables on dependent variables (responses) in bioprocesses. To Introduce to the student experimental design tools used to plan, perform full and fractional factorial designs, analyze the results, model the process based on the empirical data and suggest optimization conditions. To familiarize the student with a commercial software on the subject.'
$ws.Rows(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Rows(12).EntireRow.AutoFit()

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = '1. The role of statistics in engineering 2. Fundamentals of applied statistics3. Analysis of Variance 4. Multiple Comparison Tests 5. Design of Experiments'
$ws.Range("C14").Value = '1. The role of statistics in engineering 2. Fundamentals of applied statistics3. Analysis of Variance 4. Multiple Comparison Tests 5. Design of Experiments'
$ws.Rows(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2019'
$ws.Range("C15").Value = '01/01/2019'
$ws.Rows(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. The role of statistics in Engineering: methods of data collection2. Fundamentals of applied statistics3. Analysis of variance: analysis of a model4. Multiple comparison tests (Tukey, Hsu)5. Design of Experiments: advantages of factorial designs in comparison to “one variable at a time” experiments; full factorial design (2 ^ k), and fractionated (2^k-p), and response surface mythology'
$ws.Range("C16").Value = '1. The role of statistics in Engineering: methods of data collection2. Fundamentals of applied statistics3. Analysis of variance: analysis of a model4. Multiple comparison tests (Tukey, Hsu)5. Design of Experiments: advantages of factorial designs in comparison to “one variable at a time” experiments; full factorial design (2 ^ k), and fractionated (2^k-p), and response surface mythology'
$ws.Rows(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17:C17").Clear()
$ws.Rows(17).EntireRow.AutoFit()

# Row 18
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5817181 - Valdeir Arantes'
$ws.Range("C18").Value = '5817181 - Valdeir Arantes'
$ws.Rows(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Rows(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'MF≥ 5,0 para aprovação 5,0'
$ws.Range("C20").Value = 'MF≥ 5,0 para aprovação 5,0'
$ws.Rows(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Rows(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22:C22").Clear()
$ws.Rows(22).EntireRow.AutoFit()

# Row 23
$ws.Range("B23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Range("A23").Clear()
$ws.Rows(23).RowHeight = 30

# Row 24
$ws.Range("B24").Value = 'LOT2023 -  Processos Bioquímicos Industriais  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOT2023 -  Processos Bioquímicos Industriais  (Requisito fraco)
'
$ws.Rows(24).RowHeight = 30

